$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '29.772.11'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.54%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.869.49'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.31%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.7287'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.71%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '241.24'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3136'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.72%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07135'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '24.39'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.57%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.08133'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.09%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.7415'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.64%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.880.75'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.335'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.15%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '92.34'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.43%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '29.780.57'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.993'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.08%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '247.06'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.46%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.35'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.95%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000007804'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.31%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.128.46'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.38%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.39%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '7.735'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.25%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.1535'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.65%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.199'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.31%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '163.69'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.98%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.52'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.007'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '

$ws.Range("E30").Value = '  -2.06%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.514'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.94%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.521'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.79%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.177'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.60%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.05295'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.69%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.228'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.87%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7383'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.41%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("E38").Value = '  +0.30%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01933'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.737'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.57%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.4463'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.45%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.8713'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.95%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.956'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.56%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '71.14'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.90%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.043.31'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.16%  '

$ws.Range("E46").Value = '  -0.17%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '103.74'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("E48").Value = '  -1.43%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.439'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.04%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.503'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.78%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.024.04'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '

